# Bizagi simulation_metrics-2.xlsx edit
# Commit message: "Modified to include the simulation metrics in the json file"
#
# The underlying semantic change is a text fix in cell A3 of the
# "simulation_metrics" worksheet: the gateway label
#   "Complete  /Accurate?"  (double space)
# becomes
#   "Complete /Accurate?"   (single space)
#
# Re-typing the corrected text causes Excel to drop the old shared-string
# entry and append the corrected one at the end of the shared-strings
# table, which is exactly what the target diff shows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("simulation_metrics")

# Fix the typo: collapse the double space between "Complete" and "/Accurate?"
$ws.Range("A3").Value = "Complete /Accurate?"

# Match the author's final cursor position/selection (A3) recorded in the
# saved sheetView.
[void]$ws.Range("A3").Select()

# The author's session also left column D/E sized to fit their header
# text ("Available Resources" / "Avg Time"); reproduce that autosizing.
$ws.Columns.Item(4).ColumnWidth = 18.14
$ws.Columns.Item(5).ColumnWidth = 7.91
